$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Flow Type" worksheet — re-label/re-order the LULC rows, insert the new
#    Net Growth / Transfer rows, and extend the table with 6 new blank rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Flow Type")

# Rows 26-29 keep the left-aligned / shaded "LULC" style (s="2") that they
# already had — only the text changes.
$ws.Cells.Item(26, 1).Value = "LULC: Harvest Live"
$ws.Cells.Item(27, 1).Value = "LULC: Harvest DOM"
$ws.Cells.Item(28, 1).Value = "LULC: Mortality Live"
$ws.Cells.Item(29, 1).Value = "LULC: Transfer DOM"

# Rows 30-35 switch to the plain/default style (no explicit style id) and
# hold the new CO2/CO/CH4 emission breakdown rows.
$lulcEmission = @(
    "LULC: Emission Live CO2",
    "LULC: Emission Live CO",
    "LULC: Emission Live CH4",
    "LULC: Emission DOM CO2",
    "LULC: Emission DOM CO",
    "LULC: Emission DOM CH4"
)
for ($i = 0; $i -lt $lulcEmission.Length; $i++) {
    $cell = $ws.Cells.Item(30 + $i, 1)
    $cell.Style = "Normal"
    $cell.Value = $lulcEmission[$i]
}

# Rows 36-40 go back to the shaded style and now hold the "Net Growth" rows
# that used to live at the bottom of the sheet.
$netGrowth = @(
    "Net Growth: Atmosphere -> Coarse Roots",
    "Net Growth: Atmosphere -> Fine Roots",
    "Net Growth: Atmosphere -> Foliage",
    "Net Growth: Atmosphere -> Merchantable",
    "Net Growth: Atmosphere -> Other Wood"
)
for ($i = 0; $i -lt $netGrowth.Length; $i++) {
    $ws.Cells.Item(36 + $i, 1).Value = $netGrowth[$i]
}

# Rows 41-43 now hold the "Transfer" rows that used to live at the bottom.
$ws.Cells.Item(41, 1).Value = "Transfer: AG Slow -> BG Slow"
$ws.Cells.Item(42, 1).Value = "Transfer: Snag Branch -> AG Fast"
$ws.Cells.Item(43, 1).Value = "Transfer: Snag Stem -> AG Medium"

# Rows 44-49: six new, empty rows (default text style id 1), extending the
# sheet's used range to A1:B49.
$blankRows = $ws.Range("A44:A49")
$blankRows.Style = "Normal"
$blankRows.NumberFormat = "@"

# Column A widens to fit the new (longer) labels; column B stops being
# hidden and becomes a normal best-fit column instead.
$ws.Columns.Item(1).ColumnWidth = 48.65
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(2).ColumnWidth = 14.8

# Re-apply the AutoFilter over the new A1:B43 range (toggle off first so the
# engine actually rewrites the reference instead of removing the filter).
$ws.AutoFilterMode = $false
$ws.Range("A1:B43").AutoFilter()

# Selection / active cell moves to A42 (matches the saved cursor position).
$ws.Activate()
$ws.Range("A42").Select()

# ---------------------------------------------------------------------------
# 2. Workbook-level defined name: _FilterDatabase now spans the bigger range.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Flow Type!_FilterDatabase") {
        $n.RefersTo = "='Flow Type'!`$A`$1:`$B`$43"
    }
}
